$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price text that looks numeric (e.g. "1.001", "26.917.95").
# Force the cell to Text format before assigning so Excel keeps the exact
# string instead of silently re-parsing it as a number, then restore the
# default "Normal" style so we do not leave a stray number format behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.917.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.880.07'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.45%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '278.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.42%  '

$ws.Range("E6").Value = '  +0.02%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5308'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3465'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.51%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06981'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.77%  '

$ws.Range("E10").Value = '  +1.04%  '

$ws.Range("E11").Value = '  -2.83%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.874.40'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.12%  '

$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.85'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.59%  '

$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.197'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.28%  '

$ws.Range("E16").Value = '  +3.64%  '

$ws.Range("E17").Value = '  +0.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008068'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.57%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.956.70'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.110.76'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.760'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.66%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.203'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.372'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '147.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.67%  '

$ws.Range("E27").Value = '  +1.32%  '

$ws.Range("E28").Value = '  -0.31%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.67%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.382'
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.341'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.35%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08905'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.75%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04960'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.96%  '

$ws.Range("E34").Value = '  +4.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7313'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.24%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.886'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.304'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.83%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.399'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01859'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.41%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5176'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.96%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9630'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.83%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '116.41'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.202'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.145'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.00%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9997'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4528'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1351'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.414'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.46%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '36.37'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.92%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05956'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.09%  '

$ws.Range("E51").Value = '  +0.00%  '
